$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename sheet to reflect new "through" date
$ws.Name = "Through 2022-10-16"

# Update title cell text
$ws.Range("B1").Value = "October 2022 (through October 16)"

# Rows 3 and 4 swap rank (re-sorted by updated totals): update neighborhood labels
$ws.Range("A3").Value = "Humboldt Park"
$ws.Range("A4").Value = "Englewood"

# Apply updated neighborhood monthly carjacking counts
$ws.Range("AP2").Value = 5
$ws.Range("L2").Value = 10
$ws.Range("V2").Value = 11
$ws.Range("AA3").Value = 5
$ws.Range("AC3").Value = 2
$ws.Range("AD3").Value = 1
$ws.Range("AF3").Value = ""
$ws.Range("AH3").Value = 2
$ws.Range("AJ3").Value = 2
$ws.Range("AL3").Value = ""
$ws.Range("AN3").Value = 1
$ws.Range("AR3").Value = 4
$ws.Range("AT3").Value = 3
$ws.Range("AU3").Value = 5
$ws.Range("AV3").Value = 4
$ws.Range("AW3").Value = 2
$ws.Range("AX3").Value = 4
$ws.Range("AY3").Value = 6
$ws.Range("B3").Value = 3
$ws.Range("BA3").Value = 4
$ws.Range("BB3").Value = 5
$ws.Range("BD3").Value = 4
$ws.Range("BE3").Value = 2
$ws.Range("BF3").Value = 2
$ws.Range("BG3").Value = 3
$ws.Range("BH3").Value = 4
$ws.Range("BI3").Value = 5
$ws.Range("BJ3").Value = 1
$ws.Range("BK3").Value = 3
$ws.Range("BL3").Value = 4
$ws.Range("BM3").Value = 1
$ws.Range("BN3").Value = 3
$ws.Range("BO3").Value = 3
$ws.Range("BP3").Value = 1
$ws.Range("BQ3").Value = ""
$ws.Range("BR3").Value = ""
$ws.Range("BS3").Value = 4
$ws.Range("BU3").Value = ""
$ws.Range("BV3").Value = 2
$ws.Range("BW3").Value = 2
$ws.Range("BX3").Value = 1
$ws.Range("BZ3").Value = ""
$ws.Range("CA3").Value = 4
$ws.Range("CC3").Value = ""
$ws.Range("D3").Value = 7
$ws.Range("E3").Value = 7
$ws.Range("F3").Value = 2
$ws.Range("G3").Value = 5
$ws.Range("H3").Value = 9
$ws.Range("I3").Value = 6
$ws.Range("J3").Value = 3
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 5
$ws.Range("M3").Value = 6
$ws.Range("N3").Value = 9
$ws.Range("O3").Value = 4
$ws.Range("P3").Value = 8
$ws.Range("Q3").Value = 4
$ws.Range("S3").Value = 3
$ws.Range("T3").Value = 4
$ws.Range("U3").Value = 7
$ws.Range("V3").Value = ""
$ws.Range("W3").Value = 7
$ws.Range("X3").Value = 7
$ws.Range("Y3").Value = 3
$ws.Range("Z3").Value = 5
$ws.Range("AA4").Value = 7
$ws.Range("AC4").Value = 4
$ws.Range("AD4").Value = 3
$ws.Range("AF4").Value = 1
$ws.Range("AH4").Value = 3
$ws.Range("AJ4").Value = 4
$ws.Range("AL4").Value = 1
$ws.Range("AN4").Value = 2
$ws.Range("AR4").Value = 3
$ws.Range("AT4").Value = 1
$ws.Range("AU4").Value = 1
$ws.Range("AV4").Value = 1
$ws.Range("AW4").Value = 7
$ws.Range("AX4").Value = 1
$ws.Range("AY4").Value = 2
$ws.Range("B4").Value = 5
$ws.Range("BA4").Value = 1
$ws.Range("BB4").Value = 2
$ws.Range("BD4").Value = 1
$ws.Range("BE4").Value = ""
$ws.Range("BF4").Value = 4
$ws.Range("BG4").Value = 1
$ws.Range("BH4").Value = 3
$ws.Range("BI4").Value = ""
$ws.Range("BJ4").Value = 3
$ws.Range("BK4").Value = 4
$ws.Range("BL4").Value = 1
$ws.Range("BM4").Value = 5
$ws.Range("BN4").Value = 4
$ws.Range("BO4").Value = 2
$ws.Range("BP4").Value = 4
$ws.Range("BQ4").Value = 4
$ws.Range("BR4").Value = 2
$ws.Range("BS4").Value = 5
$ws.Range("BU4").Value = 5
$ws.Range("BV4").Value = 3
$ws.Range("BW4").Value = 3
$ws.Range("BX4").Value = 2
$ws.Range("BZ4").Value = 2
$ws.Range("C4").Value = 6
$ws.Range("CA4").Value = 1
$ws.Range("CC4").Value = 3
$ws.Range("D4").Value = 8
$ws.Range("E4").Value = 8
$ws.Range("F4").Value = 9
$ws.Range("G4").Value = 13
$ws.Range("H4").Value = 10
$ws.Range("I4").Value = 8
$ws.Range("J4").Value = 4
$ws.Range("K4").Value = 13
$ws.Range("L4").Value = 4
$ws.Range("M4").Value = 2
$ws.Range("N4").Value = 2
$ws.Range("O4").Value = 6
$ws.Range("P4").Value = 7
$ws.Range("Q4").Value = 7
$ws.Range("S4").Value = 2
$ws.Range("T4").Value = ""
$ws.Range("U4").Value = 4
$ws.Range("V4").Value = 5
$ws.Range("W4").Value = 2
$ws.Range("X4").Value = 2
$ws.Range("Y4").Value = 9
$ws.Range("Z4").Value = 8
$ws.Range("L6").Value = 8
$ws.Range("V6").Value = 8
$ws.Range("L8").Value = 2
$ws.Range("B9").Value = 3
$ws.Range("V17").Value = 4
$ws.Range("BJ20").Value = 2
$ws.Range("B24").Value = 1
$ws.Range("B32").Value = 2
$ws.Range("V32").Value = 1
$ws.Range("AP37").Value = 1
$ws.Range("BJ45").Value = 2
$ws.Range("AF48").Value = 1
$ws.Range("L51").Value = 2
$ws.Range("B65").Value = 4
$ws.Range("AP71").Value = 2
$ws.Range("V75").Value = 1
$ws.Range("AZ97").Value = 1
